$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.136.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.29%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.861.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.85%  "
# Row 4
$ws.Range("E4").Value = "  -0.30%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.06%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.35%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4670"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.79%  "
# Row 8
$ws.Range("E8").Value = "  -1.60%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06537"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.22%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.15%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07817"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.07%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "95.97"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.66%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.856.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.94%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.118"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.06%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6694"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.34%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "276.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.16%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.161.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.26%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9999"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.43%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.464"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.28%  "
# Row 20
$ws.Range("E20").Value = "  -1.47%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.095.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.12%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007238"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.91%  "
# Row 23
$ws.Range("E23").Value = "  -0.34%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.143"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.02%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.305"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.07%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.76%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.88%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.906"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.41%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.347"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.98%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09542"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.48%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.392"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.44%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.466"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.91%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.089"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.47%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04650"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.81%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7004"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.31%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.097"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.99%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.716"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.48%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01870"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.99%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.300"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.63%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.530"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.26%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.09%  "
# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.926"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.80%  "
# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8442"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.56%  "
# Row 44
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4162"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.25%  "
# Row 45
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9998"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.44%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.14%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "993.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.64%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.134"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.99%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.156"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.62%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.07%  "
# Row 51
$ws.Range("E51").Value = "  -4.53%  "
